$d = $word.ActiveDocument

$section = $d.Sections.Item(1)
$footer = $section.Footers.Item(1)
$footer.PageNumbers.Add(2)
